$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author removed a row that was sorted into the "Interesting Cohorts"
# block (andreas / d) so that only cohorts belonging to analysis 2 remain.
# That row was row 95 in the original sheet; deleting it shifts every
# following row up by one (174 data rows -> 173).
$ws.Rows(95).Delete()

# Land the selection roughly where Excel would leave it after deleting the
# row while D94 was part of the active selection.
$ws.Range("D94").Select()
